# Update payroll and system: add extra "null-numeric / empty-allowed / required"
# notes to the (OPR)Import System Ethnicity workbook's description sheet, and
# make that description sheet the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- description (รายละเอียด) sheet updates ---------------------------------

# Row 1 header note gains a clarifying suffix.
$ws2.Range("A1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# New two-column (merged) note block in columns F:G, matching the existing
# default (unstyled) cell look except for center alignment.
$ws2.Range("F1:G1").ColumnWidth = 24.8
$ws2.Range("F1:G1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("F1:G1").Merge()

$ws2.Range("F2").Value = "เป็นค่าว่างได้"
$ws2.Range("G2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# Move the live selection to C10 and make this sheet the active tab (this also
# clears tabSelected on Sheet1 and flips workbookView's activeTab to 1).
$ws2.Range("C10").Select()
